# Alteraçoes em alguns asserts e adiciona report e 1 cenario
# (changes scoped to this workbook: fix two assert values on Planilha1
#  and update the active sheet/selection to reflect where work left off)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")

# Row 3 ("Testa_Erro") asserts: fix the password casing and replace the
# placeholder numeric confirmation with the real (lowercase) password text.
$ws1.Range("D3").Value = "96Miguel."
$ws1.Range("E3").Value = "96miguel."

# Leave the workbook focused on Planilha1 with F4 selected (previously
# Planilha2/C3 was the active view).
$ws1.Activate()
$ws1.Range("F4").Select()
